# Auto-generated edit script: updates computed market-price columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) across several
# worksheets, reflecting refreshed pricing data pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 779.8
$ws.Range("I29").Value = 779.8
$ws.Range("K29").Value = 2339.4
$ws.Range("M29").Value = -2058.4
$ws.Range("H31").Value = 230
$ws.Range("I31").Value = 230
$ws.Range("K31").Value = 690
$ws.Range("M31").Value = -460
$ws.Range("H41").Value = 539.6667
$ws.Range("J41").Value = 266
$ws.Range("L41").Value = 266
$ws.Range("N41").Value = -1146
$ws.Range("H112").Value = 3178.6428
$ws.Range("J112").Value = 3158.75
$ws.Range("L112").Value = 9476.25
$ws.Range("N112").Value = -11692.25
$ws.Range("H113").Value = 6212.3
$ws.Range("I113").Value = 4037.6667
$ws.Range("K113").Value = 4037.6667
$ws.Range("M113").Value = -783.6667000000002
$ws.Range("H129").Value = 4741.25
$ws.Range("J129").Value = 4924.857
$ws.Range("L129").Value = 14774.571
$ws.Range("N129").Value = -24774.571
$ws.Range("H137").Value = 3166.6667
$ws.Range("I137").Value = 3500
$ws.Range("J137").Value = 2500
$ws.Range("K137").Value = 10500
$ws.Range("L137").Value = 7500
$ws.Range("M137").Value = -7950
$ws.Range("N137").Value = -12600

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H27").Value = 4499.5
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 4499.5
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 4499.5
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -4867.5
$ws.Range("H122").Value = 1238
$ws.Range("I122").Value = 1238
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3714
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1264
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 401
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H20").Value = 69990
$ws.Range("J20").Value = 69990
$ws.Range("L20").Value = 69990
$ws.Range("N20").Value = -70462
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H30").Value = 69990
$ws.Range("J30").Value = 69990
$ws.Range("L30").Value = 69990
$ws.Range("N30").Value = -70172
$ws.Range("H31").Value = 4383
$ws.Range("I31").Value = 4136.4
$ws.Range("J31").Value = 4999.5
$ws.Range("K31").Value = 4136.4
$ws.Range("L31").Value = 4999.5
$ws.Range("M31").Value = -3841.4
$ws.Range("N31").Value = -5589.5
$ws.Range("H34").Value = 4383
$ws.Range("I34").Value = 4136.4
$ws.Range("J34").Value = 4999.5
$ws.Range("K34").Value = 4136.4
$ws.Range("L34").Value = 4999.5
$ws.Range("M34").Value = -3934.4
$ws.Range("N34").Value = -5403.5
$ws.Range("H58").Value = 4796.2856
$ws.Range("I58").Value = 3266.6667
$ws.Range("J58").Value = 5943.5
$ws.Range("K58").Value = 3266.6667
$ws.Range("L58").Value = 5943.5
$ws.Range("M58").Value = -3063.6667
$ws.Range("N58").Value = -6349.5
$ws.Range("H128").Value = 69990
$ws.Range("J128").Value = 69990
$ws.Range("L128").Value = 69990
$ws.Range("N128").Value = -79950
$ws.Range("H134").Value = 4178.4287
$ws.Range("I134").Value = 4041.5
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 12124.5
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -9589.5
$ws.Range("N134").Value = -20070
$ws.Range("H136").Value = 4796.2856
$ws.Range("I136").Value = 3266.6667
$ws.Range("J136").Value = 5943.5
$ws.Range("K136").Value = 9800.000100000001
$ws.Range("L136").Value = 17830.5
$ws.Range("M136").Value = -7250.000100000001
$ws.Range("N136").Value = -22930.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1002
$ws.Range("J22").Value = 1002
$ws.Range("L22").Value = 3006
$ws.Range("N22").Value = -3344
$ws.Range("H25").Value = 780.2
$ws.Range("I25").Value = 475.25
$ws.Range("J25").Value = 2000
$ws.Range("K25").Value = 1425.75
$ws.Range("L25").Value = 6000
$ws.Range("M25").Value = -1256.75
$ws.Range("N25").Value = -6338
$ws.Range("H27").Value = 1002
$ws.Range("J27").Value = 1002
$ws.Range("L27").Value = 3006
$ws.Range("N27").Value = -3210
$ws.Range("H29").Value = 960
$ws.Range("J29").Value = 800
$ws.Range("L29").Value = 2400
$ws.Range("N29").Value = -2954
$ws.Range("H30").Value = 780.2
$ws.Range("I30").Value = 475.25
$ws.Range("J30").Value = 2000
$ws.Range("K30").Value = 1425.75
$ws.Range("L30").Value = 6000
$ws.Range("M30").Value = -1323.75
$ws.Range("N30").Value = -6204
$ws.Range("H131").Value = 2106.1177
$ws.Range("I131").Value = 1709.5454
$ws.Range("K131").Value = 5128.6362
$ws.Range("M131").Value = -88.63619999999992
$ws.Range("H138").Value = 11711.25
$ws.Range("I138").Value = 11711.25
$ws.Range("K138").Value = 35133.75
$ws.Range("M138").Value = -29993.75
$ws.Range("H139").Value = 3712.5
$ws.Range("I139").Value = 4900
$ws.Range("K139").Value = 14700
$ws.Range("M139").Value = -9560

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 42690.2
$ws.Range("J15").Value = 42690.2
$ws.Range("L15").Value = 42690.2
$ws.Range("N15").Value = -43266.2
$ws.Range("H81").Value = 42690.2
$ws.Range("J81").Value = 42690.2
$ws.Range("L81").Value = 42690.2
$ws.Range("N81").Value = -44686.2
$ws.Range("H84").Value = 42690.2
$ws.Range("J84").Value = 42690.2
$ws.Range("L84").Value = 128070.6
$ws.Range("N84").Value = -138054.6
$ws.Range("H107").Value = 1962.1818
$ws.Range("I107").Value = 1058.7273
$ws.Range("J107").Value = 2865.6365
$ws.Range("K107").Value = 1058.7273
$ws.Range("L107").Value = 2865.6365
$ws.Range("M107").Value = 861.2727
$ws.Range("N107").Value = -6705.636500000001
$ws.Range("H132").Value = 3648.8076
$ws.Range("I132").Value = 3428.1428
$ws.Range("J132").Value = 4575.6
$ws.Range("K132").Value = 10284.4284
$ws.Range("L132").Value = 13726.8
$ws.Range("M132").Value = -7754.428400000001
$ws.Range("N132").Value = -18786.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1695
$ws.Range("I22").Value = 1256.1428
$ws.Range("K22").Value = 1256.1428
$ws.Range("M22").Value = -961.1428000000001
$ws.Range("H27").Value = 1695
$ws.Range("I27").Value = 1256.1428
$ws.Range("K27").Value = 1256.1428
$ws.Range("M27").Value = -1149.1428
$ws.Range("H82").Value = 2270.1428
$ws.Range("I82").Value = 1974.5
$ws.Range("J82").Value = 2664.3333
$ws.Range("K82").Value = 1974.5
$ws.Range("L82").Value = 2664.3333
$ws.Range("M82").Value = -1613.5
$ws.Range("N82").Value = -3386.3333
$ws.Range("H85").Value = 2270.1428
$ws.Range("I85").Value = 1974.5
$ws.Range("J85").Value = 2664.3333
$ws.Range("K85").Value = 1974.5
$ws.Range("L85").Value = 2664.3333
$ws.Range("M85").Value = -726.5
$ws.Range("N85").Value = -5160.3333
$ws.Range("H100").Value = 2928.8333
$ws.Range("I100").Value = 2928.8333
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2928.8333
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -2387.8333
$ws.Range("N100").ClearContents()
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H132").Value = 2380.7222
$ws.Range("I132").Value = 2309.6875
$ws.Range("J132").Value = 2949
$ws.Range("K132").Value = 6929.0625
$ws.Range("L132").Value = 8847
$ws.Range("M132").Value = -4399.0625
$ws.Range("N132").Value = -13907

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3963.7144
$ws.Range("J81").Value = 6224
$ws.Range("L81").Value = 12448
$ws.Range("N81").Value = -14570
$ws.Range("H84").Value = 3963.7144
$ws.Range("J84").Value = 6224
$ws.Range("L84").Value = 62240
$ws.Range("N84").Value = -72848
$ws.Range("H122").Value = 7749.75
$ws.Range("I122").Value = 7749.75
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 23249.25
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -20799.25
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 2786.3635
$ws.Range("I132").Value = 1164.5
$ws.Range("K132").Value = 3493.5
$ws.Range("M132").Value = -963.5

